# Correct the estimated and actual hours for several activities
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 55: Estimated/Actual hours 5.0 -> 3.0
$ws.Range("D55").Value = 3.0
$ws.Range("E55").Value = 3.0

# Row 56: Estimated/Actual hours 6.0 -> 2.0, Status Pending -> Done, Comments cleared
$ws.Range("D56").Value = 2.0
$ws.Range("E56").Value = 2.0
$ws.Range("F56").Value = "Done"
$ws.Range("G56").Value = ""

# Row 57: Estimated hours 6.0 -> 3.0, Actual hours 6.0 -> 4.0
$ws.Range("D57").Value = 3.0
$ws.Range("E57").Value = 4.0

# Row 60: Estimated hours 6.0 -> 3.0, Actual hours 6.0 -> 4.0
$ws.Range("D60").Value = 3.0
$ws.Range("E60").Value = 4.0

# Row 61: Estimated hours 6.0 -> 3.0, Actual hours 6.0 -> 4.0
$ws.Range("D61").Value = 3.0
$ws.Range("E61").Value = 4.0
